# Apply edits described by the diff:
# 1. Rename sheet 1 "Basic Clinic Data..." -> "Test Import Survey Response 1"
# 2. Rename sheet 2 "Facility Fundamentals" -> "Test Import Survey Response 2"
# 3. Make sheet 2 the active tab (activeTab=1 on workbookView, tabSelected on sheet2 sheetView)
# 4. Remove tabSelected from sheet1's sheetView (it is no longer selected)
# 5. Update sheet2's selection from A6 to H24

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Test Import Survey Response 1"
$ws2.Name = "Test Import Survey Response 2"

# Update selection on sheet 2 before activating it, then activate it so it
# becomes the selected/active tab (tabSelected + workbook activeTab).
$ws2.Select()
$ws2.Range("H24").Select()

$wb.Save()
